$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.95"
$ws.Range("E2").Value = "'1.02%"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'32.05"
$ws.Range("E3").Value = "'1.17%"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.011"
$ws.Range("E4").Value = "'-1.57%"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07874"
$ws.Range("E5").Value = "'-3.53%"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.119"
$ws.Range("E6").Value = "'-16.99%"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.806"
$ws.Range("E7").Value = "'0.52%"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.784"
$ws.Range("E8").Value = "'-1.60%"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9273"
$ws.Range("E9").Value = "'-0.44%"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1745"
$ws.Range("E10").Value = "'-0.91%"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07902"
$ws.Range("E11").Value = "'4.94%"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08808"
$ws.Range("E12").Value = "'-1.91%"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03115"
$ws.Range("E13").Value = "'3.85%"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.02%"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001515"
$ws.Range("E15").Value = "'0.61%"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005983"
$ws.Range("E16").Value = "'2.99%"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.451"
$ws.Range("E17").Value = "'-3.64%"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.270"
$ws.Range("E18").Value = "'0.58%"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'1.38%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1289"
$ws.Range("E20").Value = "'-2.92%"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.157"
$ws.Range("E21").Value = "'6.32%"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'5.50%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'-0.26%"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001235"
$ws.Range("E24").Value = "'-0.42%"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004498"
$ws.Range("E25").Value = "'0.81%"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001247"
$ws.Range("E26").Value = "'4.08%"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.01735"
$ws.Range("E39").Value = "'-1.70%"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04763"
$ws.Range("E40").Value = "'4.86%"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007419"
$ws.Range("E41").Value = "'7.73%"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1364"
$ws.Range("E42").Value = "'0.79%"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002135"
$ws.Range("E43").Value = "'-3.25%"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01074"
$ws.Range("E44").Value = "'8.27%"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006060"
$ws.Range("E45").Value = "'-6.21%"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.003394"
$ws.Range("E47").Value = "'-61.17%"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.8205"
$ws.Range("E48").Value = "'-0.01%"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002096"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0001996"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
